$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert a brand-new first worksheet "1. D&A" before the current first
#    sheet, and rename/renumber the existing sheets to make room for it.
#    NOTE: worksheet variables in this runtime track *position*, not object
#    identity, so every sheet reference is (re)fetched by index right after
#    any operation that could change the sheet ordering.
# ---------------------------------------------------------------------------
$newSheet = $wb.Worksheets.Add($wb.Worksheets.Item(1))
$newSheet.Name = "1. D&A"

$gapSheet    = $wb.Worksheets.Item(2)   # was "1. Trade de GAP em Cripto"
$binarySheet = $wb.Worksheets.Item(3)   # was "2. Opção binária"
$socialSheet = $wb.Worksheets.Item(4)   # was "3. Rede social profissional"

$gapSheet.Name    = "2. Trade de GAP em Cripto"
$binarySheet.Name = "3. Opção binária"
$socialSheet.Name = "4. Rede social profissional"

# Re-fetch the new sheet by position (it stays first).
$ds = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# 2. Populate the new sheet with the "D&A" (Definir & Alinhar) checklist.
# ---------------------------------------------------------------------------
function Set-HeaderCell($cell, $val) {
    $cell.Value = $val
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4131   # xlLeft
    $cell.VerticalAlignment = -4108     # xlCenter
    $cell.WrapText = $true
}

function Set-DataCell($cell, $val) {
    $cell.Value = $val
    $cell.HorizontalAlignment = -4131   # xlLeft
    $cell.VerticalAlignment = -4108     # xlCenter
    $cell.WrapText = $true
}

Set-HeaderCell $ds.Range("A1") "Setor"
Set-HeaderCell $ds.Range("B1") "Item"
Set-HeaderCell $ds.Range("C1") "Definição"
Set-HeaderCell $ds.Range("D1") "Prioridade"
Set-HeaderCell $ds.Range("E1") "Status"

$rows = @(
    @("Estratégico", "Definir problema", 1, "Em andamento"),
    @("Estratégico", "Definir propósito", 2, "Em andamento"),
    @("Estratégico", "Definir público alvo", 3, "Em andamento"),
    @("Estratégico", "Definir escopo de atuação (iniciar com MVP) - Só dashboard, PaaS, SaaS, Consultoria etc", 4, "Em andamento"),
    @("Estratégico", "Sonhar alto: Pensar em escopo grandioso daqui X anos - (para ter um norte)", 5, "Em andamento"),
    @("Estratégico", "Definir nome da marca", 6, "Em andamento"),
    @("Estratégico", "Definir slogan/branding inicial", 7, "Em andamento"),
    @("Estratégico", "Definir funções e responsabilidades", 8, "Em andamento"),
    @("Estratégico", "Definir pré-requisitos para início de piloto", 9, "Em andamento"),
    @("Estratégico", "Definir metas/deadlines para entrega dos pré-requisitos", 10, "Em andamento"),
    @("Estratégico", "Definir metas/deadlines para piloto do 1º MVP", 11, "Em andamento")
)

$r = 2
foreach ($row in $rows) {
    Set-DataCell $ds.Cells.Item($r, 1) $row[0]
    Set-DataCell $ds.Cells.Item($r, 2) $row[1]
    Set-DataCell $ds.Cells.Item($r, 4) $row[2]
    Set-DataCell $ds.Cells.Item($r, 5) $row[3]
    $r++
}

# ---------------------------------------------------------------------------
# 3. Match the look & feel (column widths, frozen header row) of the other
#    sheets in the workbook.
# ---------------------------------------------------------------------------
$ds.Columns.Item(1).ColumnWidth = 16.85546875
$ds.Columns.Item(2).ColumnWidth = 81
$ds.Columns.Item(3).ColumnWidth = 125.140625
$ds.Columns.Item(4).ColumnWidth = 12.5703125
$ds.Columns.Item(5).ColumnWidth = 14.42578125

$ds.Application.ActiveWindow.SplitRow = 1
$ds.Application.ActiveWindow.FreezePanes = $true

# ---------------------------------------------------------------------------
# 4. Turn the populated range into a real table ("Table135"), mirroring the
#    other three sheets which each hold their own table.
# ---------------------------------------------------------------------------
$tableRange = $ds.Range("A1:E1048575")
$lo = $ds.ListObjects.Add(1, $tableRange, 0, 1)
$lo.Name = "Table135"
$lo.TableStyle = "TableStyleMedium2"

Write-Host "D&A sheet created with $($rows.Count) rows and table Table135"
